$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("41:41").Insert()
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value = 45246
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 100112010
$ws.Range("G41").Value = "Achicoria"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 55
$ws.Range("K41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = 10000
$ws.Range("N41").Value = "$/caja 18 unidades"
$ws.Range("O41").Value = "Región Metropolitana"
$ws.Range("P41").Value = 556
$ws.Range("Q41").Value = 18
$ws.Range("R41").Value = "Hortaliza"
